$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated e-mail message: add Dawson's e-mail address next to
# "Anne-Sophie Hollender" (row 2) in the E-mail column (F).
$ws.Range("F2").Value = "DawsonCoding@gmail.com"

# Added completion percentage / tax status: "Bianca Saia" (row 12)
# now shows as having paid ("Y") instead of not ("N").
$ws.Range("C12").Value = "Y"

# Leave the selection where the user finished editing.
$ws.Range("F14").Select()
